$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "StatQuery" (column C) used by all three tabs (Cases/Samples/Files):
# replaces the old, much longer filter-driven StatQuery with a simpler
# Programs/Studies/Cases/Samples/Case Files/Study Files count query.
$statQuery = @'
MATCH (p:program)<--(s:study)<-[*]-(c:case)<--(demo:demographic)
OPTIONAL MATCH (samp:sample)-->(c)
OPTIONAL MATCH (diag:diagnosis)-->(c)
OPTIONAL MATCH (f:file)-[*]->(c)
OPTIONAL MATCH (sf:file)-->(s)
WITH DISTINCT f, sf, samp AS samp, c, demo, diag, s, p
 WHERE samp.specific_sample_pathology IN ["Osteosarcoma"]  
RETURN  
    count(distinct p) AS Programs,
    count(distinct s) AS Studies,
    count(distinct c) AS Cases,
    count(distinct samp) AS Samples,
    count(distinct f) AS `Case Files`,
    count(distinct sf) AS `Study Files`
'@

# FilesTab "query" (column B, row 4): same as before but the trailing
# "Study Code" output column (and its MATCH) has been dropped.
$filesQuery = @'
MATCH (f:file)-->(parent)
WITH DISTINCT f, parent
MATCH (f)-[*]->(c:case)<--(demo:demographic)
 MATCH (s:study)<-[*]-(c)<--(diag:diagnosis)
 MATCH (samp:sample)-->(c) 
 WHERE samp.specific_sample_pathology IN ["Osteosarcoma"]  
WITH DISTINCT f, parent, c, demo, diag, s
RETURN coalesce(f.file_name, '') AS `File Name`, 
        coalesce(labels(parent)[0], '') AS `Association`,
        coalesce(f.file_description, '') AS `Description`,
        coalesce(f.file_format, '') AS `Format`,
        coalesce(f.file_size, '') AS `Size`,
        coalesce(c.case_id, '') AS `Case ID`, 
        coalesce(diag.disease_term,'') AS Diagnosis 
'@

# CasesTab "query" (column B, row 2): same as before but now also returns
# a `Cohort` column sourced from co.cohort_description.
$casesQuery = @'
MATCH (s:study)<-[*]-(c:case)<--(demo:demographic)
MATCH (c)<--(diag:diagnosis)
MATCH (samp:sample)-->(c) 
 WHERE samp.specific_sample_pathology IN ["Osteosarcoma"]  
OPTIONAL MATCH (co:cohort)<-[*]-(c)
  WITH DISTINCT c, s, demo, diag, co
RETURN  coalesce(c.case_id, '') AS `Case ID` ,
        coalesce(s.clinical_study_designation, '') AS `Study Code` ,
        coalesce(s.clinical_study_type, '') AS  `Study Type`,
        coalesce(demo.breed, '') AS Breed ,
        coalesce(diag.disease_term, '') AS Diagnosis ,
        coalesce(diag.stage_of_disease, '') AS `Stage of Disease` ,
        coalesce(demo.patient_age_at_enrollment, '') AS Age ,
        coalesce(demo.sex, '') AS Sex ,
        coalesce(demo.neutered_indicator, '') AS `Neutered Status`,
        coalesce(demo.weight, '') AS `Weight (kg)`,
        coalesce(diag.best_response, '') AS `Response to Treatment`,
coalesce(co.cohort_description, '') AS `Cohort`
'@

# Here-strings always carry a trailing newline from the closing "'@" line;
# the source cells do not, so trim exactly that one trailing newline.
$statQuery  = $statQuery.TrimEnd("`r","`n")
$filesQuery = $filesQuery.TrimEnd("`r","`n")
$casesQuery = $casesQuery.TrimEnd("`r","`n")

# CasesTab row (row 2)
$ws.Range("B2").Value2 = $casesQuery
$ws.Range("C2").Value2 = $statQuery

# SamplesTab row (row 3) - query (B3) is unchanged, only StatQuery (C3) updates
$ws.Range("C3").Value2 = $statQuery

# FilesTab row (row 4)
$ws.Range("B4").Value2 = $filesQuery
$ws.Range("C4").Value2 = $statQuery

# Row heights shrink (they were all pegged at the 409.6 max before)
$ws.Rows.Item(2).RowHeight = 270
$ws.Rows.Item(3).RowHeight = 225
$ws.Rows.Item(4).RowHeight = 210

# Zoom level goes from 25% to 70%
$excel.ActiveWindow.Zoom = 70
